# Add the "London_Enviornment" worksheet (terrain / landclass per station)
# and select it as the active sheet, matching the authored workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the selection on the existing neighbor sheet before we leave it.
$null = $ws1.Range("B2:B25").Select()

# New sheet is inserted right after "London_Nearest_Neighbor".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "London_Enviornment"

# header + 34 rows: 24 met/aq stations (same order as sheet 1) followed by
# the 10 "london_grid_*" rows sorted alphabetically.
$data = @(
    @("station_id","terrain","landclass"),
    @("BX9","flat","suburbs"),
    @("BX1","flat","suburbs"),
    @("BL0","flat","park"),
    @("CD9","flat","city"),
    @("CD1","flat","suburbs"),
    @("CT2","flat","city"),
    @("CT3","flat","city"),
    @("CR8","flat","park"),
    @("GN0","flat","suburbs"),
    @("GR4","flat","park"),
    @("GN3","flat","park"),
    @("GR9","flat","suburbs"),
    @("GB0","flat","suburbs"),
    @("HR1","flat","park"),
    @("HV1","flat","suburbs"),
    @("LH0","flat","farm"),
    @("KC1","flat","suburbs"),
    @("KF1","flat","suburbs"),
    @("LW2","flat","suburbs"),
    @("RB7","flat","industrial"),
    @("TD5","flat","industrial"),
    @("ST5","flat","park"),
    @("TH4","flat","city"),
    @("MY7","flat","city"),
    @("london_grid_346","flat","suburbs"),
    @("london_grid_366","flat","park"),
    @("london_grid_368","flat","suburbs"),
    @("london_grid_388","flat","park"),
    @("london_grid_408","flat","suburbs"),
    @("london_grid_409","flat","suburbs"),
    @("london_grid_430","flat","industrial"),
    @("london_grid_451","flat","park"),
    @("london_grid_452","flat","park"),
    @("london_grid_472","flat","farm")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# The bottom block (grid cells) is kept alphabetically sorted by station_id.
$ws2.Range("A26:C35").Sort($ws2.Range("A26"), 1)

# size the columns to fit their content and leave the cursor where the
# author left it, one row below the last data row.
$ws2.Columns("A:C").AutoFit()
$null = $ws2.Range("B36").Select()
$ws2.Activate()
